# Scheduled runner update: refresh market-board derived values (profit
# sheets) for several Leve rows across the ALC/ARM/BSM/CRP/CUL/LTW/WVR
# worksheets. Cells are plain numeric values (no formulas in this
# workbook), so we just overwrite the ones that changed.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 883
$ws.Range("I19").Value = 825
$ws.Range("K19").Value = 825
$ws.Range("M19").Value = -650
$ws.Range("H43").Value = 3000
$ws.Range("J43").Value = 3000
$ws.Range("L43").Value = 3000
$ws.Range("N43").Value = -3138
$ws.Range("H55").Value = 392.5
$ws.Range("I55").Value = 390
$ws.Range("K55").Value = 390
$ws.Range("M55").Value = -176
$ws.Range("H86").Value = 2247.5
$ws.Range("J86").Value = 2495
$ws.Range("L86").Value = 2495
$ws.Range("N86").Value = -4741
$ws.Range("H87").Value = 95925
$ws.Range("J87").Value = 95925
$ws.Range("L87").Value = 95925
$ws.Range("N87").Value = -98421
$ws.Range("H89").Value = 2247.5
$ws.Range("J89").Value = 2495
$ws.Range("L89").Value = 12475
$ws.Range("N89").Value = -23707
$ws.Range("H90").Value = 95925
$ws.Range("J90").Value = 95925
$ws.Range("L90").Value = 287775
$ws.Range("N90").Value = -300255
$ws.Range("H98").Value = 476.42856
$ws.Range("I98").Value = 476.42856
$ws.Range("K98").Value = 476.42856
$ws.Range("M98").Value = 1021.57144
$ws.Range("H122").Value = 476.42856
$ws.Range("I122").Value = 476.42856
$ws.Range("K122").Value = 1429.28568
$ws.Range("M122").Value = 1020.71432
$ws.Range("H132").Value = 1194.8462
$ws.Range("I132").Value = 1194.8462
$ws.Range("K132").Value = 3584.5386
$ws.Range("M132").Value = -1054.5386
$ws.Range("H138").Value = 5410.1377
$ws.Range("J138").Value = 5736.077
$ws.Range("L138").Value = 17208.231
$ws.Range("N138").Value = -27488.231

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 996.6667
$ws.Range("I2").Value = 996.6667
$ws.Range("K2").Value = 996.6667
$ws.Range("M2").Value = -883.6667
$ws.Range("H32").Value = 7496.4165
$ws.Range("I32").Value = 7496.4165
$ws.Range("K32").Value = 7496.4165
$ws.Range("M32").Value = -7209.4165
$ws.Range("H61").Value = 2216.75
$ws.Range("I61").Value = 1024.0714
$ws.Range("J61").Value = 4999.6665
$ws.Range("K61").Value = 1024.0714
$ws.Range("L61").Value = 4999.6665
$ws.Range("M61").Value = -812.0714
$ws.Range("N61").Value = -5423.6665
$ws.Range("H63").Value = 7302
$ws.Range("J63").Value = 20006
$ws.Range("L63").Value = 20006
$ws.Range("N63").Value = -21378
$ws.Range("H66").Value = 7302
$ws.Range("J66").Value = 20006
$ws.Range("L66").Value = 100030
$ws.Range("N66").Value = -106894
$ws.Range("H74").Value = 1513.3334
$ws.Range("I74").Value = 1513.3334
$ws.Range("K74").Value = 1513.3334
$ws.Range("M74").Value = -639.3334
$ws.Range("H77").Value = 1513.3334
$ws.Range("I77").Value = 1513.3334
$ws.Range("K77").Value = 7566.666999999999
$ws.Range("M77").Value = -3198.666999999999
$ws.Range("H102").Value = 2646.8333
$ws.Range("I102").Value = 2764.3333
$ws.Range("J102").Value = 2529.3333
$ws.Range("K102").Value = 2764.3333
$ws.Range("L102").Value = 2529.3333
$ws.Range("M102").Value = -1142.3333
$ws.Range("N102").Value = -5773.3333
$ws.Range("H116").Value = 996.6667
$ws.Range("I116").Value = 996.6667
$ws.Range("K116").Value = 996.6667
$ws.Range("M116").Value = 1297.3333
$ws.Range("H136").Value = 2216.75
$ws.Range("I136").Value = 1024.0714
$ws.Range("J136").Value = 4999.6665
$ws.Range("K136").Value = 3072.2142
$ws.Range("L136").Value = 14998.9995
$ws.Range("M136").Value = -522.2142000000003
$ws.Range("N136").Value = -20098.9995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 996.6667
$ws.Range("I3").Value = 996.6667
$ws.Range("K3").Value = 996.6667
$ws.Range("M3").Value = -882.6667
$ws.Range("H107").Value = 799.3333
$ws.Range("I107").Value = 799.3333
$ws.Range("K107").Value = 799.3333
$ws.Range("M107").Value = 1120.6667
$ws.Range("H134").Value = 18850
$ws.Range("I134").Value = 15315.833
$ws.Range("K134").Value = 45947.499
$ws.Range("M134").Value = -43412.499

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 287.25
$ws.Range("I22").Value = 287.5
$ws.Range("J22").Value = 287
$ws.Range("K22").Value = 287.5
$ws.Range("L22").Value = 287
$ws.Range("M22").Value = 62.5
$ws.Range("N22").Value = -987
$ws.Range("H31").Value = 4444.5835
$ws.Range("I31").Value = 2778.3333
$ws.Range("K31").Value = 2778.3333
$ws.Range("M31").Value = -2483.3333
$ws.Range("H34").Value = 4444.5835
$ws.Range("I34").Value = 2778.3333
$ws.Range("K34").Value = 2778.3333
$ws.Range("M34").Value = -2576.3333
$ws.Range("H58").Value = 995
$ws.Range("I58").Value = 995
$ws.Range("K58").Value = 995
$ws.Range("M58").Value = -792
$ws.Range("H132").Value = 2230.6924
$ws.Range("I132").Value = 1772.7273
$ws.Range("J132").Value = 4749.5
$ws.Range("K132").Value = 5318.1819
$ws.Range("L132").Value = 14248.5
$ws.Range("M132").Value = -2788.1819
$ws.Range("N132").Value = -19308.5
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 300
$ws.Range("I134").Value = 300
$ws.Range("K134").Value = 900
$ws.Range("M134").Value = 1635
$ws.Range("H136").Value = 995
$ws.Range("I136").Value = 995
$ws.Range("K136").Value = 2985
$ws.Range("M136").Value = -435
$ws.Range("H141").Value = 99986
$ws.Range("J141").Value = 99986
$ws.Range("L141").Value = 99986
$ws.Range("N141").Value = -110346

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3660.36
$ws.Range("I34").Value = 3242
$ws.Range("K34").Value = 9726
$ws.Range("M34").Value = -9642
$ws.Range("H55").Value = 9399.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2081.8333
$ws.Range("J22").Value = 2249
$ws.Range("L22").Value = 2249
$ws.Range("N22").Value = -2839
$ws.Range("H27").Value = 2081.8333
$ws.Range("J27").Value = 2249
$ws.Range("L27").Value = 2249
$ws.Range("N27").Value = -2463
$ws.Range("H132").Value = 5499.5
$ws.Range("I132").Value = 6499.5
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 19498.5
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -16968.5
$ws.Range("N132").Value = -18558.5
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 50000
$ws.Range("N134").Value = -60140
$ws.Range("H136").Value = 2130
$ws.Range("I136").Value = 1945
$ws.Range("K136").Value = 5835
$ws.Range("M136").Value = -3285

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1514.2858
$ws.Range("I113").Value = 1250
$ws.Range("J113").Value = 1620
$ws.Range("K113").Value = 3750
$ws.Range("L113").Value = 4860
$ws.Range("M113").Value = -1580
$ws.Range("N113").Value = -9200
$ws.Range("H132").Value = 2870.5881
$ws.Range("I132").Value = 1927.1923
$ws.Range("K132").Value = 5781.5769
$ws.Range("M132").Value = -3251.5769
$ws.Range("H136").Value = 1248.1578
$ws.Range("I136").Value = 1133.4615
$ws.Range("K136").Value = 3400.3845
$ws.Range("M136").Value = -850.3845
